# Regenerate the localization-status report for a new handoff run.
# New source file id: 09c2564d-9e3f-4e64-b5f3-f06f832c5af7
# New xliff content hash: 200d9454f0c6994c646b0d2f93ae82418e5e23f9

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"
$overview.Range("B2").Value = "e2e\09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"
$overview.Range("G2").Value = "2016-09-05 03:05:35"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"
$zhcn.Range("G2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.200d9454f0c6994c646b0d2f93ae82418e5e23f9.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-09-05 03:05:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.md"
$dede.Range("G2").Value = "09c2564d-9e3f-4e64-b5f3-f06f832c5af7.200d9454f0c6994c646b0d2f93ae82418e5e23f9.de-de.xlf"
$dede.Range("H2").Value = "2016-09-05 03:05:35"
